# Auto-generated Excel COM-interop script to revert CareSocialCodes term sheet
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update Metadata sheet (sheet1) ---
$ws1.Cells.Item(2,2).Value = 'http://fhir.kl.dk/term/CodeSystem/CareSocialCodes'
$ws1.Cells.Item(3,2).Value = '1.1.0'
$ws1.Cells.Item(4,2).Value = 'CareSocialCodes'
$ws1.Cells.Item(5,2).Value = 'CareSocialCodes'
$ws1.Cells.Item(6,2).Value = 'active'
$ws1.Cells.Item(7,2).Value = 'false'
$ws1.Cells.Item(8,2).Value = '2023-07-10T23:08:03+02:00'
$ws1.Cells.Item(9,2).Value = 'Kommunernes Landsforening'
$ws1.Cells.Item(10,2).Value = 'No display for ContactDetail'
$ws1.Cells.Item(11,2).Value = 'Administrative/technical codes in Local Govenment Denmark (KL), associated with KLCommonCareSocial'
$ws1.Cells.Item(12,2).Value = ''
$ws1.Cells.Item(13,2).Value = ''
$ws1.Cells.Item(14,2).Value = 'true'
$ws1.Cells.Item(15,2).Value = ''
$ws1.Cells.Item(16,2).Value = 'is-a'
$ws1.Cells.Item(17,2).Value = ''
$ws1.Cells.Item(18,2).Value = ''
$ws1.Cells.Item(19,2).Value = 'complete'
$ws1.Cells.Item(20,2).Value = ''
$ws1.Cells.Item(21,2).Value = ''

# --- Rewrite Concepts sheet (sheet2) ---
# Remove existing data rows (2..75), keep header row 1
for ($r = 75; $r -ge 2; $r--) {
    $ws2.Rows.Item($r).Delete()
}

# Write target data rows
$ws2.Cells.Item(2,1).Value = '1'
$ws2.Cells.Item(2,2).Value = '95ec4535-8fe8-4296-867c-35de421794cf'
$ws2.Cells.Item(2,3).Value = 'evaluering'
$ws2.Cells.Item(2,4).Value = 'evaluering'
$ws2.Cells.Item(3,1).Value = '2'
$ws2.Cells.Item(3,2).Value = 'effe55c7-572c-4a99-8fb4-2a9dda2f6572'
$ws2.Cells.Item(3,3).Value = 'FFB støttebehovsvurdering'
$ws2.Cells.Item(3,4).Value = 'Angiver, hvor stort et behov borgeren har for hjælp og støtte.'
$ws2.Cells.Item(4,1).Value = '2'
$ws2.Cells.Item(4,2).Value = '053a301d-1bb8-4cc4-b781-87825ecf0ef8'
$ws2.Cells.Item(4,3).Value = 'FFB vurdering af borgers situation'
$ws2.Cells.Item(4,4).Value = 'Vurdering der sammenholder oplysninger fra udredning mhp en samlet faglig analyse og konklusion.'
$ws2.Cells.Item(5,1).Value = '2'
$ws2.Cells.Item(5,2).Value = '54c4ffea-7caf-4edc-8aa9-ef6e0be26c4c'
$ws2.Cells.Item(5,3).Value = 'FSIII helhedsvurdering'
$ws2.Cells.Item(5,4).Value = 'Sagsvurdering, der sammenholder resultaterne af en myndighedsudredning og den øvrige sagsoplysning med henblik på at træffe en afgørelse'
$ws2.Cells.Item(6,1).Value = '2'
$ws2.Cells.Item(6,2).Value = '3f7a8ca0-afca-4b0d-8773-a99b5f2f8aaf'
$ws2.Cells.Item(6,3).Value = 'VUM Borgerens perspektiv på indsatsen'
$ws2.Cells.Item(6,4).Value = 'VUM Borgerens perspektiv på indsatsen'
$ws2.Cells.Item(7,1).Value = '2'
$ws2.Cells.Item(7,2).Value = 'f52887de-023f-4193-b6b0-4b0a37b1cffc'
$ws2.Cells.Item(7,3).Value = 'VUM Borgerens ressourcer i forhold til indsatsen'
$ws2.Cells.Item(7,4).Value = 'VUM Borgerens ressourcer i forhold til indsatsen'
$ws2.Cells.Item(8,1).Value = '1'
$ws2.Cells.Item(8,2).Value = '7b41185e-eeb4-437d-8120-5d51bbd27a79'
$ws2.Cells.Item(8,3).Value = 'Indsats/ydelses-anmodning'
$ws2.Cells.Item(8,4).Value = 'Indsats/ydelses-anmodning'
$ws2.Cells.Item(9,1).Value = '2'
$ws2.Cells.Item(9,2).Value = 'ad865929-7363-4b2d-a271-01993181fbaf'
$ws2.Cells.Item(9,3).Value = 'Hjemmepleje'
$ws2.Cells.Item(9,4).Value = 'Der anmodes om hjemmepleje'
$ws2.Cells.Item(10,1).Value = '2'
$ws2.Cells.Item(10,2).Value = '490ab7be-ddb1-4a54-baf1-009fe6e8a83b'
$ws2.Cells.Item(10,3).Value = 'Hjemmesygepleje'
$ws2.Cells.Item(10,4).Value = 'Der anmodes om sygepleje'
$ws2.Cells.Item(11,1).Value = '2'
$ws2.Cells.Item(11,2).Value = '4a297733-4d66-4726-a933-590d55cf91e0'
$ws2.Cells.Item(11,3).Value = 'Social indsats'
$ws2.Cells.Item(11,4).Value = 'Der anmodes om en social indsats'
$ws2.Cells.Item(12,1).Value = '2'
$ws2.Cells.Item(12,2).Value = 'a71921ea-fe83-441d-933b-cc21d0b3c8c3'
$ws2.Cells.Item(12,3).Value = 'Sundhedsfremme og forebyggelse'
$ws2.Cells.Item(12,4).Value = 'Der anmodes om sundhedsfremme og forebyggelse'
$ws2.Cells.Item(13,1).Value = '2'
$ws2.Cells.Item(13,2).Value = '7fc66c15-0cb3-4c89-9e18-f3a01e6a6592'
$ws2.Cells.Item(13,3).Value = 'Træning'
$ws2.Cells.Item(13,4).Value = 'Der anmodes om træning'
$ws2.Cells.Item(14,1).Value = '1'
$ws2.Cells.Item(14,2).Value = 'b5731132-f6b9-4a47-995d-681d2b755d4f'
$ws2.Cells.Item(14,3).Value = 'Kilde'
$ws2.Cells.Item(14,4).Value = 'Person der er kilde til oplysning eller vurdering'
$ws2.Cells.Item(15,1).Value = '2'
$ws2.Cells.Item(15,2).Value = '25b4e705-2e9a-47a2-b11a-c829316b9d3a'
$ws2.Cells.Item(15,3).Value = 'Borger'
$ws2.Cells.Item(15,4).Value = 'Borger er kilde til oplysning eller vurdering'
$ws2.Cells.Item(16,1).Value = '2'
$ws2.Cells.Item(16,2).Value = 'f6ea2920-7dde-491e-a489-6b99a3904069'
$ws2.Cells.Item(16,3).Value = 'Sagsbehandler'
$ws2.Cells.Item(16,4).Value = 'Sagsbehandler er kilde til oplysning eller vurdering'
$ws2.Cells.Item(17,1).Value = '2'
$ws2.Cells.Item(17,2).Value = '63338442-7b2e-405b-acc0-142361ef19f1'
$ws2.Cells.Item(17,3).Value = 'Andre'
$ws2.Cells.Item(17,4).Value = 'Andre end borger og sagsbehandler er kilde til oplysning eller vurdering. Fx læge, pårørende og udfører'
$ws2.Cells.Item(18,1).Value = '3'
$ws2.Cells.Item(18,2).Value = 'f00a6844-1005-401d-965d-1c5859df7beb'
$ws2.Cells.Item(18,3).Value = 'Udfører'
$ws2.Cells.Item(18,4).Value = 'Udfører er kilde til oplysning eller vurdering'
$ws2.Cells.Item(19,1).Value = '1'
$ws2.Cells.Item(19,2).Value = '25303acd-dcaf-4a8e-a8a3-3961a43858aa'
$ws2.Cells.Item(19,3).Value = 'kontaktaktivitet'
$ws2.Cells.Item(19,4).Value = 'aktivitet der foretages ved kommunal kontakt.'
$ws2.Cells.Item(20,1).Value = '2'
$ws2.Cells.Item(20,2).Value = 'bb6fc544-7f4f-4b50-8868-1431e0df2381'
$ws2.Cells.Item(20,3).Value = 'observation'
$ws2.Cells.Item(20,4).Value = 'observation foretaget ved kontakt'
$ws2.Cells.Item(21,1).Value = '2'
$ws2.Cells.Item(21,2).Value = '9f03dfbb-7a97-45a5-94db-d4c3501714a9'
$ws2.Cells.Item(21,3).Value = 'opfølgning'
$ws2.Cells.Item(21,4).Value = 'opfølgning foretaget ved kontakt'
$ws2.Cells.Item(22,1).Value = '2'
$ws2.Cells.Item(22,2).Value = '829ac647-c7fc-4964-836b-f708d886e0e3'
$ws2.Cells.Item(22,3).Value = 'oplysning'
$ws2.Cells.Item(22,4).Value = 'oplysning er foretaget ved kontakt'
$ws2.Cells.Item(23,1).Value = '2'
$ws2.Cells.Item(23,2).Value = '15775b0a-7ec6-469e-9d3c-a81fbc9a1b45'
$ws2.Cells.Item(23,3).Value = 'udførelse af akutindsats'
$ws2.Cells.Item(23,4).Value = 'Akutindsats, bevilget af læge eller akutteam, udført ved kontakt'
$ws2.Cells.Item(24,1).Value = '2'
$ws2.Cells.Item(24,2).Value = 'c03b426a-4348-407f-b343-f4baa9759c72'
$ws2.Cells.Item(24,3).Value = 'udførelse af ikke-bevilget indsats'
$ws2.Cells.Item(24,4).Value = 'Indsats, der ikke i forvejen var bevilget, udført ved kontakt'
$ws2.Cells.Item(25,1).Value = '2'
$ws2.Cells.Item(25,2).Value = '784275f1-6822-4a88-b361-d958007d5253'
$ws2.Cells.Item(25,3).Value = 'udførelse af planlagt indsats'
$ws2.Cells.Item(25,4).Value = 'Planlagt indsats udført ved kontakt'
$ws2.Cells.Item(26,1).Value = '2'
$ws2.Cells.Item(26,2).Value = '9269c9a2-8220-447b-a127-811275b41062'
$ws2.Cells.Item(26,3).Value = 'vurdering/bevilling'
$ws2.Cells.Item(26,4).Value = 'vurdering og/eller bevilling er foretaget ved kontakt'
$ws2.Cells.Item(27,1).Value = '1'
$ws2.Cells.Item(27,2).Value = '3762dd32-4123-43a8-815d-ec40d3697652'
$ws2.Cells.Item(27,3).Value = 'indforståelse ifm henvendelse'
$ws2.Cells.Item(27,4).Value = 'indforståelse ifm henvendelse'
$ws2.Cells.Item(28,1).Value = '2'
$ws2.Cells.Item(28,2).Value = '54aeeba6-6aa2-4165-a5a8-bbd6f2f3b1eb'
$ws2.Cells.Item(28,3).Value = 'ikke indforstået'
$ws2.Cells.Item(28,4).Value = 'Borger ikke indforstået med henvendelse'
$ws2.Cells.Item(29,1).Value = '2'
$ws2.Cells.Item(29,2).Value = 'e67035da-9179-466b-99ad-ea86835d38c9'
$ws2.Cells.Item(29,3).Value = 'indforstået'
$ws2.Cells.Item(29,4).Value = 'Borger indforstået med henvendelse'
$ws2.Cells.Item(30,1).Value = '1'
$ws2.Cells.Item(30,2).Value = '253bbdc0-c4ca-4e77-9d3e-3a9e51281636'
$ws2.Cells.Item(30,3).Value = 'Mål/formål'
$ws2.Cells.Item(30,4).Value = 'Mål/formål'
$ws2.Cells.Item(31,1).Value = '2'
$ws2.Cells.Item(31,2).Value = '416fe27d-3ccf-4390-8742-8b52a9d8dc78'
$ws2.Cells.Item(31,3).Value = 'FFB borgers mål og ønsker'
$ws2.Cells.Item(31,4).Value = 'Borgers mål og ønsker som specificeret af FFB. Det er overordnet og er som udtrykt af borger'
$ws2.Cells.Item(32,1).Value = '2'
$ws2.Cells.Item(32,2).Value = '6746d4af-145a-4bfd-a672-05c0cf11b53b'
$ws2.Cells.Item(32,3).Value = 'FFB delmål'
$ws2.Cells.Item(32,4).Value = 'Delmål i FFB, ligger under indsatsmål, og er som udtrykt af udfører'
$ws2.Cells.Item(33,1).Value = '2'
$ws2.Cells.Item(33,2).Value = '0bb3daef-538d-45dc-b444-abdbcb63f6bc'
$ws2.Cells.Item(33,3).Value = 'FFB indsatsmål'
$ws2.Cells.Item(33,4).Value = 'Indsatsmål i FFB - er de konkrete faglige mål der arbejdes med, og som er knyttet til tilstande og opfølgninger'
$ws2.Cells.Item(34,1).Value = '2'
$ws2.Cells.Item(34,2).Value = 'ffb9886b-d04e-46b1-9165-a400f91f822b'
$ws2.Cells.Item(34,3).Value = 'FSIII borgers ønsker og mål'
$ws2.Cells.Item(34,4).Value = 'FSIII borgers ønsker og mål'
$ws2.Cells.Item(35,1).Value = '2'
$ws2.Cells.Item(35,2).Value = 'ca552020-6ed1-4cdc-b0d4-32697f1f27ad'
$ws2.Cells.Item(35,3).Value = 'FSIII tilstandsmål'
$ws2.Cells.Item(35,4).Value = 'Tilstandsmål som defineret af FSIII - ofte i form af en forventet tilstand'
$ws2.Cells.Item(36,1).Value = '2'
$ws2.Cells.Item(36,2).Value = '424827b1-23aa-4848-962b-56ee47def560'
$ws2.Cells.Item(36,3).Value = 'Indsatsformål'
$ws2.Cells.Item(36,4).Value = 'Det overordnede formål med hele indsatsen (Defineret af FFB men bruges også for FSIII)'
$ws2.Cells.Item(37,1).Value = '1'
$ws2.Cells.Item(37,2).Value = '11266a8f-5795-42ab-88ec-4fe5c6c28e80'
$ws2.Cells.Item(37,3).Value = 'målfokus'
$ws2.Cells.Item(37,4).Value = 'målfokus'
$ws2.Cells.Item(38,1).Value = '2'
$ws2.Cells.Item(38,2).Value = '66959f77-6e2a-4574-8423-3ff097f8b9fa'
$ws2.Cells.Item(38,3).Value = 'funktionsevneniveau'
$ws2.Cells.Item(38,4).Value = 'målet udtrykkes i form af et funktionsevneniveau, som kan være et FFB funktionsevneniveau eller FSIII tilstandsniveau'
$ws2.Cells.Item(39,1).Value = '2'
$ws2.Cells.Item(39,2).Value = '90c48f03-f194-4b2f-ad7d-6cba1069ae48'
$ws2.Cells.Item(39,3).Value = 'måltype'
$ws2.Cells.Item(39,4).Value = 'målet udtrykkes i form af en måltype, måltypen er den forventede ændring i tilstanden givet indsatsen'
$ws2.Cells.Item(40,1).Value = '1'
$ws2.Cells.Item(40,2).Value = '940f37e6-8a3d-483b-adac-be8af3268a5b'
$ws2.Cells.Item(40,3).Value = 'oplysningsaktivitet'
$ws2.Cells.Item(40,4).Value = 'oplysningsaktivitet'
$ws2.Cells.Item(41,1).Value = '2'
$ws2.Cells.Item(41,2).Value = '95e787e0-5490-437d-ae4c-f3736644242f'
$ws2.Cells.Item(41,3).Value = 'afklarende samtale §119, FSIII'
$ws2.Cells.Item(41,4).Value = 'udførelse af afklarende samtale vedr. sundhedsfremme og forebyggelse jævnfør FSIII'
$ws2.Cells.Item(42,1).Value = '2'
$ws2.Cells.Item(42,2).Value = 'e5a73b0e-a5d2-430e-931f-6156306ab00f'
$ws2.Cells.Item(42,3).Value = 'funktionsevnevurdering hjemmepleje, FSIII'
$ws2.Cells.Item(42,4).Value = 'udførelse af funktionsevnevurdering i hjemmeplejen jævnfør FSIII'
$ws2.Cells.Item(43,1).Value = '2'
$ws2.Cells.Item(43,2).Value = 'e70c66c0-a939-493b-8ea8-5b7e7b48ba1a'
$ws2.Cells.Item(43,3).Value = 'generelle oplysninger, FSIII'
$ws2.Cells.Item(43,4).Value = 'indhentning af generelle oplysninger jævnfør FSIII'
$ws2.Cells.Item(44,1).Value = '2'
$ws2.Cells.Item(44,2).Value = 'f8ebd11a-04f3-4aa0-9786-406e8896c84d'
$ws2.Cells.Item(44,3).Value = 'socialfaglig udredning, VUM/FFB'
$ws2.Cells.Item(44,4).Value = 'udførelse af socialfaglig udredning som specificeret af VUM og FFB'
$ws2.Cells.Item(45,1).Value = '2'
$ws2.Cells.Item(45,2).Value = '47fd1468-89da-4803-9d7a-ecc039a30d92'
$ws2.Cells.Item(45,3).Value = 'sygeplejefaglig udredning, FSIII'
$ws2.Cells.Item(45,4).Value = 'udførelse af sygeplejefaglig udredning jævnfør FSIII'
$ws2.Cells.Item(46,1).Value = '2'
$ws2.Cells.Item(46,2).Value = '0f0f223c-abe3-4720-aab8-c257679a0a4e'
$ws2.Cells.Item(46,3).Value = 'terapeutfaglig udredning, FSIII'
$ws2.Cells.Item(46,4).Value = 'udførelse af terapeutfaglig udredning jævnfør FSIII'
$ws2.Cells.Item(47,1).Value = '1'
$ws2.Cells.Item(47,2).Value = '10deb210-e7e1-4d56-9531-b9ff2102126e'
$ws2.Cells.Item(47,3).Value = 'Planlagt indsatsforløb'
$ws2.Cells.Item(47,4).Value = 'Planlagt indsatsforløb'
$ws2.Cells.Item(48,1).Value = '2'
$ws2.Cells.Item(48,2).Value = 'e459386d-1474-4c31-89c5-c8bc7a25e3d4'
$ws2.Cells.Item(48,3).Value = 'Social indsats'
$ws2.Cells.Item(48,4).Value = 'Det planlagte indsatsforløb er en social indsats, som defineret af FFB, med tilhørende ydelser, målgruppe og tilbud'
$ws2.Cells.Item(49,1).Value = '2'
$ws2.Cells.Item(49,2).Value = '4fd6c23a-6ff3-4251-ac37-3ca095027b5b'
$ws2.Cells.Item(49,3).Value = 'Sundhedsfremme og forebyggelse §119'
$ws2.Cells.Item(49,4).Value = 'Det planlagte indsatsforløb vedrører sundhedsfremme og forebyggelse efter §119'
$ws2.Cells.Item(50,1).Value = '3'
$ws2.Cells.Item(50,2).Value = '5c160c02-e858-4c1f-925a-71ed64844749'
$ws2.Cells.Item(50,3).Value = 'Interventionsforløb efter §119'
$ws2.Cells.Item(50,4).Value = 'Et planlagt §119 indsatsforløb, som udelukkende indeholder interventioner, og som leveres samlet. Med interventioner menes indsatser, der søger at forbedre borgers tilstand. For §119 kunne et interventionsforløb fx være et KOL-træningshold der indeholder indsatserne fysisk træning og funktionstræning.'
$ws2.Cells.Item(51,1).Value = '3'
$ws2.Cells.Item(51,2).Value = '9791e55a-656f-47eb-8fd5-c4a06b0a4662'
$ws2.Cells.Item(51,3).Value = 'Opfølgningsforløb efter §119'
$ws2.Cells.Item(51,4).Value = 'Planlagt indsatsforløb der består af opfølgningsindsatser. Anvendes hvis der i forvejen planlægges opfølgning ved flere forskellige lejligheder fx efter 2,4 og 6 måneder.'
$ws2.Cells.Item(52,1).Value = '2'
$ws2.Cells.Item(52,2).Value = 'ddd2f670-5ec7-4f9c-9a2c-aee25cb133bf'
$ws2.Cells.Item(52,3).Value = 'Genoptræning efter §140'
$ws2.Cells.Item(52,4).Value = 'Det planlagte indsatsforløb er genoptræning efter sundhedslovens §140'
$ws2.Cells.Item(53,1).Value = '3'
$ws2.Cells.Item(53,2).Value = 'f15b2663-94d9-4d0c-a5de-d8bd8e1e4ebb'
$ws2.Cells.Item(53,3).Value = 'Interventionsforløb efter §140'
$ws2.Cells.Item(53,4).Value = 'Et planlagt §140 indsatsforløb, som udelukkende indeholder interventioner, og som leveres samlet. Med interventioner menes indsatser, der søger at forbedre borgers tilstand. For §140 kunne et interventionsforløb fx være et Knæ-genoptræningshold der indeholder indsatserne fysisk træning, funktionstræning og vejledning og undervisning.'
$ws2.Cells.Item(54,1).Value = '3'
$ws2.Cells.Item(54,2).Value = '4863001e-14c7-4be8-a2da-e4f21a4b6ac4'
$ws2.Cells.Item(54,3).Value = 'Opfølgningsforløb efter §140'
$ws2.Cells.Item(54,4).Value = 'Planlagt indsatsforløb efter §140 der består af opfølgningsindsatser. Anvendes hvis der i forvejen planlægges opfølgning ved flere forskellige lejligheder fx efter 2,4 og 6 måneder.'
$ws2.Cells.Item(55,1).Value = '1'
$ws2.Cells.Item(55,2).Value = '2c059407-fed5-4852-92d8-6bb5ad63d7bb'
$ws2.Cells.Item(55,3).Value = 'Begrundelse for indsatsophør'
$ws2.Cells.Item(55,4).Value = 'Begrundelse for, at kommunale indsatser ophører'
$ws2.Cells.Item(56,1).Value = '2'
$ws2.Cells.Item(56,2).Value = '82e99421-31da-4915-96ed-168ccfa1d20c'
$ws2.Cells.Item(56,3).Value = 'Hændelse medfører ophør'
$ws2.Cells.Item(56,4).Value = 'Hændelse, som ikke er et aktivt fravalg eller en behovsvurdering medfører ophør. Det kan fx være en længere hospitalsindlæggelse, eller at borger flytter.'
$ws2.Cells.Item(57,1).Value = '2'
$ws2.Cells.Item(57,2).Value = '4bbf6d6a-a1c6-41c2-b8c1-7352b7378adf'
$ws2.Cells.Item(57,3).Value = 'Ikke yderligere behov (borger-vurderet)'
$ws2.Cells.Item(57,4).Value = 'Borger vurderer, at han/hun ikke har yderligere behov og afslutter derfor før tid'
$ws2.Cells.Item(58,1).Value = '2'
$ws2.Cells.Item(58,2).Value = 'a63b6aa6-7d56-4e67-a5b3-d73f6d262af5'
$ws2.Cells.Item(58,3).Value = 'Ikke yderligere behov (fagperson-vurderet)'
$ws2.Cells.Item(58,4).Value = 'Fagperson vurderer, at borger ikke har yderligere behov og afslutter derfor før tid'
$ws2.Cells.Item(59,1).Value = '2'
$ws2.Cells.Item(59,2).Value = '3a5a72b7-addf-4857-b80c-04e4246e3072'
$ws2.Cells.Item(59,3).Value = 'Aktivt fravalg'
$ws2.Cells.Item(59,4).Value = 'Borger har behov, men har foretaget et aktivt fravalg'
$ws2.Cells.Item(60,1).Value = '3'
$ws2.Cells.Item(60,2).Value = '0cd5734d-b663-47c6-a3da-6b14a937d144'
$ws2.Cells.Item(60,3).Value = 'Aktivt fravalg pga. anden sygdom'
$ws2.Cells.Item(60,4).Value = 'Borger har behov, men har foretaget et aktivt fravalg pga. anden sygdom'
$ws2.Cells.Item(61,1).Value = '3'
$ws2.Cells.Item(61,2).Value = '8371b769-4bfb-4ac8-b130-d91c54784a56'
$ws2.Cells.Item(61,3).Value = 'Aktivt fravalg pga. logistik ifm. transport'
$ws2.Cells.Item(61,4).Value = 'Borger har behov, men har foretaget et aktivt fravalg pga. den logistiske udfordring det er mht. transport, at nå frem til det sted indsatsen tilbydes'
$ws2.Cells.Item(62,1).Value = '3'
$ws2.Cells.Item(62,2).Value = 'a3f2bd01-078b-486e-81be-797d192ad7bd'
$ws2.Cells.Item(62,3).Value = 'Aktivt fravalg pga. anden træning'
$ws2.Cells.Item(62,4).Value = 'Borger har behov, men har foretaget et aktivt fravalg fordi han/hun er påbegyndt træning i andet regi fx fitnesscenter'
$ws2.Cells.Item(63,1).Value = '3'
$ws2.Cells.Item(63,2).Value = 'f8436c2e-af1c-44fe-939d-473b518dd01d'
$ws2.Cells.Item(63,3).Value = 'Aktivt fravalg pga. økonomi ifm. transport'
$ws2.Cells.Item(63,4).Value = 'Borger har behov, men har foretaget et aktivt fravalg pga. den omkostning der er forbundet med transport til det sted indsatsen tilbydes'
$ws2.Cells.Item(64,1).Value = '1'
$ws2.Cells.Item(64,2).Value = '3f79cee2-b148-4f2c-9cbd-387820e74685'
$ws2.Cells.Item(64,3).Value = 'Leveringstype'
$ws2.Cells.Item(64,4).Value = 'Om der leveres individuelt eller gruppebaseret'
$ws2.Cells.Item(65,1).Value = '2'
$ws2.Cells.Item(65,2).Value = '2865f123-15a7-4a36-a514-32ea37c400ca'
$ws2.Cells.Item(65,3).Value = 'Gruppebaseret indsats'
$ws2.Cells.Item(65,4).Value = 'Indsats leveres i gruppe'
$ws2.Cells.Item(66,1).Value = '2'
$ws2.Cells.Item(66,2).Value = '8d12d74c-17da-47a7-a4fe-e69dbaec0a8c'
$ws2.Cells.Item(66,3).Value = 'Individuel indsats'
$ws2.Cells.Item(66,4).Value = 'Indsats leveres individuelt'
